# Rebuild the open/close/high/low/shares_outstanding/fixed_ticker columns
# (D:I, rows 2-44) for NORTHROP_GRUMMAN_CORP with the corrected, company-specific
# data. Previously these columns had been polluted with price/volume/ticker data
# pulled in from other companies' files; "fixed_ticker" is reset to NOC throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @(2, 134.8845040186177, 129.7870483398438, 139.8218723361719, 128.7254346349996, 142720067),
    @(3, 135.1842254137789, 146.5060272216797, 149.4613873512652, 133.8801487392348, 142720067),
    @(4, 141.2568559303741, 159.7649078369141, 165.0748087824118, 138.4317125792599, 142720067),
    @(5, 158.9351089023986, 158.1488952636719, 165.1051901296775, 149.551805015991, 142720067),
    @(6, 168.7124312880878, 177.0111694335938, 179.3540495369891, 167.914316794299, 142720067),
    @(7, 191.7541241273548, 186.6952667236328, 193.1502623061722, 185.8248238563632, 142720067),
    @(8, 184.3420841267015, 198.1893768310547, 199.0548326250768, 183.493941146186, 142720067),
    @(9, 203.336216709597, 198.9759063720703, 206.6542120713813, 194.4592567403336, 142720067),
    @(10, 207.1169090338772, 214.4224090576172, 217.944384457451, 205.0072089886013, 142720067),
    @(11, 225.3783092230785, 230.2791748046875, 233.499733695505, 224.608169769281, 142720067),
    @(12, 253.4178197543208, 259.5929565429688, 269.3255944634218, 252.2934714009841, 142720067),
    @(13, 270.8566916840437, 300.0970458984375, 303.331298951087, 267.2082431760531, 142720067),
    @(14, 308.7819250257306, 284.7247924804688, 319.0643453429059, 281.1617551693087, 142720067),
    @(15, 271.929628001382, 266.6496887207031, 288.6479191999652, 254.7587628413754, 142720067),
    @(16, 283.3748506355448, 233.391128540039, 286.6447438813564, 224.2229796754844, 142720067),
    @(17, 215.8703942910444, 246.643798828125, 254.6996663214705, 212.2094533414537, 142720067),
    @(18, 243.7046294130246, 260.5849914550781, 265.0882280991191, 237.7452635963541, 142720067),
    @(19, 293.3756624858548, 311.944091796875, 324.0311547663087, 286.9484831492365, 142720067),
    @(20, 340.6062015108164, 319.32568359375, 341.6208497043885, 308.2007435049304, 142720067),
    @(21, 313.2135399240307, 340.6119995117188, 350.1055255200535, 312.4678733494119, 142720067),
    @(22, 268.184560771225, 301.7886352539062, 325.9284228137339, 265.236665383571, 142720067),
    @(23, 283.0903419627515, 297.8832397460937, 301.9618063213848, 261.6617354875416, 142720067),
    @(24, 292.0761147463565, 266.7386169433594, 296.7883512441894, 262.4221283191555, 142720067),
    @(25, 281.1592962103258, 265.0148315429688, 285.5698816732484, 261.5658932780586, 142720067),
    @(26, 300.141566381451, 329.3157958984375, 331.5828406463527, 298.5341858912371, 142720067),
    @(27, 340.1320590455784, 338.7231140136719, 348.9495833117484, 327.1717295298723, 142720067),
    @(28, 338.9704786111626, 334.7442321777344, 382.3573377193686, 329.8526659311078, 142720067),
    @(29, 364.1859519933977, 348.1389770507812, 384.91051694712, 345.2966523595877, 142720067),
    @(30, 421.9831395971285, 415.1986083984375, 450.9733606086371, 408.5747265843274, 142720067),
    @(31, 453.8910402326041, 454.2040100097656, 462.7398967130165, 411.4961246955878, 142720067),
    @(32, 453.4574284390929, 522.5267944335938, 529.436594916117, 440.7704264249136, 142720067),
    @(33, 520.3574169372527, 427.8443603515625, 523.1076040263433, 411.5151463665732, 142720067),
    @(34, 444.0294556903269, 442.0838928222656, 463.2934207338669, 424.5738270097146, 142720067),
    @(35, 436.1446182208199, 428.3382873535156, 444.701772488369, 415.4496223149758, 142720067),
    @(36, 425.7623497286203, 455.7506408691406, 480.3638846399089, 400.7720923602251, 142720067),
    @(37, 455.3914992951702, 433.6211547851562, 469.6009331883021, 413.4716803220234, 142720067),
    @(38, 466.2285226483278, 472.69970703125, 478.1281152471265, 434.8763259579104, 142720067),
    @(39, 428.3641692447456, 474.0981140136719, 479.0513169851881, 413.7688510475956, 142720067),
    @(40, 519.076124522949, 500.2463684082031, 545.9940366958643, 493.1999694655907, 142720067),
    @(41, 465.8143382725705, 480.8944091796875, 501.0275138424837, 444.9114706757644, 142720067),
    @(42, 510.0833043890325, 482.2861938476562, 540.0415861670613, 446.2312164266363, 142720067),
    @(43, 497.9926921320322, 574.3525390625, 575.1195439429852, 487.603502808621, 142720067),
    @(44, 607.219970703125, 605.5800170898438, 640.9000244140625, 571.3599853515625, 142720067)
)

foreach ($row in $rowData) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = $row[5]
    $ws.Range("I$r").Value = "NOC"
}
